# Auto-generated edit script applying the Pandaemonium_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Range("H33").Value = 284.83334
$ws.Range("I33").Value = 188.79167
$ws.Range("K33").Value = 188.79167
$ws.Range("M33").Value = 40.20832999999999

# row 43
$ws.Range("H43").Value = 777.8570999999999
$ws.Range("J43").Value = 890.0909
$ws.Range("L43").Value = 890.0909
$ws.Range("N43").Value = -1028.0909

# row 69
$ws.Range("H69").Value = 4975.8
$ws.Range("J69").Value = 4975.8
$ws.Range("L69").Value = 14927.4
$ws.Range("N69").Value = -16675.4

# row 72
$ws.Range("H72").Value = 4975.8
$ws.Range("J72").Value = 4975.8
$ws.Range("L72").Value = 44782.2
$ws.Range("N72").Value = -53518.2

# row 74
$ws.Range("H74").Value = 5117.8335
$ws.Range("I74").Value = 4601
$ws.Range("J74").Value = 5634.6665
$ws.Range("K74").Value = 4601
$ws.Range("L74").Value = 5634.6665
$ws.Range("M74").Value = -3665
$ws.Range("N74").Value = -7506.6665

# row 77
$ws.Range("H77").Value = 5117.8335
$ws.Range("I77").Value = 4601
$ws.Range("J77").Value = 5634.6665
$ws.Range("K77").Value = 23005
$ws.Range("L77").Value = 28173.3325
$ws.Range("M77").Value = -18325
$ws.Range("N77").Value = -37533.3325

# row 97
$ws.Range("H97").Value = 1508.4286
$ws.Range("J97").Value = 1508.4286
$ws.Range("L97").Value = 4525.2858
$ws.Range("N97").Value = -5517.2858

# row 138
$ws.Range("H138").Value = 5054.6714
$ws.Range("I138").Value = 1224.0714
$ws.Range("J138").Value = 7438.156
$ws.Range("K138").Value = 3672.2142
$ws.Range("L138").Value = 22314.468
$ws.Range("M138").Value = 1467.7858
$ws.Range("N138").Value = -32594.468

# row 141
$ws.Range("H141").Value = 1560.2368
$ws.Range("I141").Value = 1230.8286
$ws.Range("J141").Value = 5403.3335
$ws.Range("K141").Value = 3692.4858
$ws.Range("L141").Value = 16210.0005
$ws.Range("M141").Value = 1487.5142
$ws.Range("N141").Value = -26570.0005

$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 7966.9062
$ws.Range("I61").Value = 6435.4443
$ws.Range("J61").Value = 9935.929
$ws.Range("K61").Value = 6435.4443
$ws.Range("L61").Value = 9935.929
$ws.Range("M61").Value = -6223.4443
$ws.Range("N61").Value = -10359.929

# row 74
$ws.Range("H74").Value = 6108.136
$ws.Range("I74").Value = 4380.5405
$ws.Range("J74").Value = 15239.714
$ws.Range("K74").Value = 4380.5405
$ws.Range("L74").Value = 15239.714
$ws.Range("M74").Value = -3506.5405
$ws.Range("N74").Value = -16987.714

# row 77
$ws.Range("H77").Value = 6108.136
$ws.Range("I77").Value = 4380.5405
$ws.Range("J77").Value = 15239.714
$ws.Range("K77").Value = 21902.7025
$ws.Range("L77").Value = 76198.57000000001
$ws.Range("M77").Value = -17534.7025
$ws.Range("N77").Value = -84934.57000000001

# row 136
$ws.Range("H136").Value = 7966.9062
$ws.Range("I136").Value = 6435.4443
$ws.Range("J136").Value = 9935.929
$ws.Range("K136").Value = 19306.3329
$ws.Range("L136").Value = 29807.787
$ws.Range("M136").Value = -16756.3329
$ws.Range("N136").Value = -34907.787

$ws = $wb.Worksheets.Item("BSM")
# row 96
$ws.Range("H96").Value = 1223
$ws.Range("I96").Value = 1223
$ws.Range("K96").Value = 1223
$ws.Range("M96").Value = 1523

$ws = $wb.Worksheets.Item("CRP")
# row 21
$ws.Range("H21").Value = 19800
$ws.Range("J21").Value = 19800
$ws.Range("L21").Value = 19800
$ws.Range("N21").Value = -20270

# row 60
$ws.Range("H60").Value = 29768.334
$ws.Range("J60").Value = 29768.334
$ws.Range("L60").Value = 29768.334
$ws.Range("N60").Value = -30790.334

$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 950.0625

# row 6
$ws.Range("H6").Value = 95230.5
$ws.Range("I6").Value = 135757.28
$ws.Range("K6").Value = 407271.84
$ws.Range("M6").Value = -407158.84

# row 17
$ws.Range("H17").Value = 1833.3334
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 1800
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 5400
$ws.Range("M17").Value = -5831
$ws.Range("N17").Value = -5738

# row 70
$ws.Range("H70").Value = 2104.818
$ws.Range("I70").Value = 1021.8571
$ws.Range("K70").Value = 3065.5713
$ws.Range("M70").Value = -2750.5713

# row 73
$ws.Range("H73").Value = 2104.818
$ws.Range("I73").Value = 1021.8571
$ws.Range("K73").Value = 3065.5713
$ws.Range("M73").Value = -1973.5713

# row 107
$ws.Range("H107").Value = 954.0345
$ws.Range("I107").Value = 338.8125
$ws.Range("J107").Value = 1711.2307
$ws.Range("K107").Value = 1016.4375
$ws.Range("L107").Value = 5133.6921
$ws.Range("M107").Value = 903.5625
$ws.Range("N107").Value = -8973.6921

# row 113
$ws.Range("H113").Value = 686.8871
$ws.Range("J113").Value = 661.48
$ws.Range("L113").Value = 1984.44
$ws.Range("N113").Value = -6324.440000000001

# row 123
$ws.Range("H123").Value = 3299.6
$ws.Range("J123").Value = 3999.5
$ws.Range("L123").Value = 11998.5
$ws.Range("N123").Value = -16898.5

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 4641.7417
$ws.Range("I80").Value = 5141.727
$ws.Range("J80").Value = 4366.75
$ws.Range("K80").Value = 5141.727
$ws.Range("L80").Value = 4366.75
$ws.Range("M80").Value = -4143.727
$ws.Range("N80").Value = -6362.75

# row 83
$ws.Range("H83").Value = 4641.7417
$ws.Range("I83").Value = 5141.727
$ws.Range("J83").Value = 4366.75
$ws.Range("K83").Value = 25708.635
$ws.Range("L83").Value = 21833.75
$ws.Range("M83").Value = -20716.635
$ws.Range("N83").Value = -31817.75

# row 122
$ws.Range("H122").Value = 6713.727
$ws.Range("I122").Value = 10805.4
$ws.Range("K122").Value = 32416.2
$ws.Range("M122").Value = -29966.2

# row 126
$ws.Range("H126").Value = 2864.353
$ws.Range("I126").Value = 1989.4
$ws.Range("J126").Value = 4114.2856
$ws.Range("K126").Value = 5968.200000000001
$ws.Range("L126").Value = 12342.8568
$ws.Range("M126").Value = -3498.200000000001
$ws.Range("N126").Value = -17282.8568

$ws = $wb.Worksheets.Item("LTW")
# row 2
$ws.Range("H2").Value = 29000
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 53000
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 53000
$ws.Range("M2").Value = -4888
$ws.Range("N2").Value = -53224

# row 35
$ws.Range("H35").Value = 25846.154
$ws.Range("J35").Value = 25846.154
$ws.Range("L35").Value = 25846.154
$ws.Range("N35").Value = -26518.154

$ws = $wb.Worksheets.Item("WVR")
# row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# row 122
$ws.Range("H122").Value = 1969.5
$ws.Range("I122").Value = 1458.9166
$ws.Range("J122").Value = 3501.25
$ws.Range("K122").Value = 4376.7498
$ws.Range("L122").Value = 10503.75
$ws.Range("M122").Value = -1926.7498
$ws.Range("N122").Value = -15403.75
